$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
WITH Distinct_Samples AS (
    SELECT DISTINCT
        sp.participant_id,
        sp.study_participant_id,
        s.study_name,
        s.phs_accession,
        sp.gender,
        smp.sample_id
    FROM 
        df_participant sp
    JOIN 
        df_study s ON sp."study.phs_accession" = s.phs_accession
    JOIN 
        df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_program p ON p.program_acronym = s."program.program_acronym"
    JOIN
        df_file f1 ON f1."sample.sample_id" = smp.sample_id
    JOIN
        df_genomic_info gi ON gi."file.file_id" = f1.file_id
    WHERE 
        s.phs_accession = 'phs003432'
),
Sample_Limit AS (
    SELECT
        participant_id,
        study_participant_id,
        study_name,
        phs_accession,
        sample_id,
        ROW_NUMBER() OVER (PARTITION BY study_participant_id ORDER BY sample_id) AS rn
    FROM 
        Distinct_Samples
)
SELECT
    f1.file_name AS "File Name",
    s.study_name AS "Study Name",
    s.phs_accession AS "Accession",
    sp.participant_id AS "Participant Id",
    COALESCE(STRING_AGG(
        CASE WHEN sl.rn <= 5 THEN sl.sample_id ELSE NULL END, ', '
    ) || CASE WHEN MAX(sl.rn) > 5 THEN ', ...' ELSE '' END, '') AS "Sample Id",
    f1.file_type AS "File Type",
    gi.library_strategy AS "Library Strategy"
FROM 
    df_study s
INNER JOIN 
    df_participant sp ON sp."study.phs_accession" = s.phs_accession
INNER JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
INNER JOIN 
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
INNER JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
LEFT JOIN 
    Sample_Limit sl ON sl.study_participant_id = sp.study_participant_id
WHERE 
    s.phs_accession = 'phs003432'
GROUP BY
    f1.file_name,
    s.study_name,
    s.phs_accession,
    sp.participant_id,
    f1.file_type,
    gi.library_strategy
ORDER BY 
    f1.file_name ASC
LIMIT 100;
'@

$ws.Range("B4").Value2 = $newQuery
$ws.Rows.Item(4).RowHeight = 409.5
$ws.Range("B4").Select()
